$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 85059
$ws.Range("B2").Value = "Emilly da Costa"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45082
$ws.Range("G2").Value = 10051.41

# Row 3
$ws.Range("A3").Value = 35302
$ws.Range("B3").Value = "Vicente Viana"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 4927.09

# Row 4
$ws.Range("A4").Value = 29892
$ws.Range("B4").Value = "Isadora Correia"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45087
$ws.Range("G4").Value = 9601.280000000001

# Row 5
$ws.Range("A5").Value = 75611
$ws.Range("B5").Value = "Laís Barbosa"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45097
$ws.Range("G5").Value = 4592.5

# Row 6
$ws.Range("A6").Value = 90264
$ws.Range("B6").Value = "Gustavo Henrique Gomes"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45097
$ws.Range("G6").Value = 12094.92

# Row 7
$ws.Range("A7").Value = 8699
$ws.Range("B7").Value = "Dra. Ana Lívia Cardoso"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Outros"
$ws.Range("F7").Value = 45088
$ws.Range("G7").Value = 5100.01

# Row 8
$ws.Range("A8").Value = 42034
$ws.Range("B8").Value = "João Guilherme Carvalho"
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 12163.74

# Row 9
$ws.Range("A9").Value = 29445
$ws.Range("B9").Value = "Sr. Benício Duarte"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45098
$ws.Range("G9").Value = 6417.65

# Row 10
$ws.Range("A10").Value = 91242
$ws.Range("B10").Value = "Fernando Moura"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45085
$ws.Range("G10").Value = 7321.45

# Row 11
$ws.Range("A11").Value = 8315
$ws.Range("B11").Value = "Cauã Cardoso"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 8419.129999999999
